$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 16.277
$ws.Range("D3").Value = -8.222
$ws.Range("E3").Value = 16.558
$ws.Range("D4").Value = -7.769
$ws.Range("E9").Value = 17.185
$ws.Range("A11").Value = -21.669
$ws.Range("A12").Value = -21.652
$ws.Range("D14").Value = -7.746
$ws.Range("A15").Value = -21.904
$ws.Range("E15").Value = 16.199
$ws.Range("E19").Value = 16.239
$ws.Range("E20").Value = 16.312
$ws.Range("E25").Value = 17.08
$ws.Range("D26").Value = -7.667999999999999
$ws.Range("A27").Value = -21.747
$ws.Range("E27").Value = 16.733
$ws.Range("A28").Value = -21.889
$ws.Range("E28").Value = 16.904
$ws.Range("E30").Value = 16.354
$ws.Range("A31").Value = -21.572
$ws.Range("D31").Value = -7.931999999999999
$ws.Range("A32").Value = -21.742
$ws.Range("E32").Value = 16.901
$ws.Range("D35").Value = -7.858
$ws.Range("A36").Value = -20.5
$ws.Range("D37").Value = -7.778
$ws.Range("A38").Value = -19.741
$ws.Range("D39").Value = -7.186
$ws.Range("D40").Value = -7.858
$ws.Range("E44").Value = 16.503
$ws.Range("D45").Value = -7.390000000000001
$ws.Range("A46").Value = -21.812
$ws.Range("E47").Value = 16.572
$ws.Range("D52").Value = -7.87
$ws.Range("A54").Value = -21.703
$ws.Range("A55").Value = -22.026
$ws.Range("A56").Value = -22.047
$ws.Range("D57").Value = -8.289999999999999
$ws.Range("E58").Value = 16.597
$ws.Range("E62").Value = 16.311
$ws.Range("A67").Value = -21.6
$ws.Range("A69").Value = -21.68
$ws.Range("A72").Value = -21.445
$ws.Range("A73").Value = -19.921
$ws.Range("E77").Value = 16.751
$ws.Range("E78").Value = 16.391
$ws.Range("D81").Value = -7.27
$ws.Range("A83").Value = -21.757
$ws.Range("D83").Value = -8.289
$ws.Range("E84").Value = 16.821
$ws.Range("A86").Value = -22.035
$ws.Range("E89").Value = 17.33900000000001
$ws.Range("A91").Value = -21.613
$ws.Range("E91").Value = 17.243
$ws.Range("E92").Value = 16.891
$ws.Range("A93").Value = -21.288
$ws.Range("E96").Value = 16.677
$ws.Range("A99").Value = -20.326
$ws.Range("D100").Value = -8.425000000000001
$ws.Range("D102").Value = -7.507
$ws.Range("E102").Value = 16.59500000000001
